# Avances Controllers y DAOs.xlsx -- "Creación y avances en Plan de Adquisiciones Pagos"
#
# 1) Bump two progress cells (Controllers!D53, Daos!C53) to 100% complete.
# 2) Add a brand-new "Avance" summary sheet at the end of the workbook that
#    rolls up progress across Controllers / Daos / Vistas plus a handful of
#    standalone milestones, with a header row, a totals row and a big
#    "Avance Total" percentage at the bottom.
# 3) Leave the UI selection/active-sheet state the way the author left it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Progress updates on existing sheets
# ---------------------------------------------------------------------
$wsControllers = $wb.Worksheets.Item("Controllers")
$wsDaos        = $wb.Worksheets.Item("Daos")

$wsControllers.Range("D53").Value = 1
$wsDaos.Range("C53").Value = 1

# ---------------------------------------------------------------------
# 2) New "Avance" sheet, placed after the last existing sheet ("Vistas")
# ---------------------------------------------------------------------
$wsVistas = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAvance = $wb.Worksheets.Add($null, $wsVistas)
$wsAvance.Name = "Avance"

# Column widths
$wsAvance.Columns.Item(2).ColumnWidth = 60.28515625
$wsAvance.Columns.Item(3).ColumnWidth = 5.42578125
$wsAvance.Columns.Item(4).ColumnWidth = 12.7109375
$wsAvance.Columns.Item(5).ColumnWidth = 11.85546875
$wsAvance.Columns.Item(6).ColumnWidth = 18.140625
$wsAvance.Columns.Item(7).ColumnWidth = 7.42578125

# --- Header row (row 2) ------------------------------------------------
$wsAvance.Range("B2").Value = "SIPRO - Migración a estandar DTI"
$wsAvance.Range("C2").Value = "Total"
$wsAvance.Range("D2").Value = "Completados"
$wsAvance.Range("E2").Value = "% de avance"
$wsAvance.Range("F2").Value = "Dias de trabajo estimado (1 persona)"
$wsAvance.Range("G2").Value = "Dias de avance"

$headerRow = $wsAvance.Range("B2:G2")
$headerRow.Borders.LineStyle = 1
$headerRow.RowHeight = 45
$wsAvance.Range("F2:G2").HorizontalAlignment = -4108
$wsAvance.Range("F2:G2").WrapText = $true

# --- Detail rows (4-11) --------------------------------------------------
$labels = @(
    "Instalación de ambiente de trabajo en computadoras personales",
    "Definición de arquitectura de software a implementar",
    "Implementación de estructura básica de arquitectura de software",
    "Programación de POCOs",
    "Programación de DAOs",
    "Programación de Controllers",
    "Programación de POCOs Analytic",
    "Programación de Vista"
)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = 4 + $i
    $wsAvance.Range("B$r").Value = $labels[$i]
}

$wsAvance.Range("C4").Value = 1
$wsAvance.Range("D4").Value = 1
$wsAvance.Range("F4").Value = 3

$wsAvance.Range("C5").Value = 1
$wsAvance.Range("D5").Value = 1
$wsAvance.Range("F5").Value = 5

$wsAvance.Range("C6").Value = 1
$wsAvance.Range("D6").Value = 1
$wsAvance.Range("F6").Value = 5

$wsAvance.Range("C7").Value = 101
$wsAvance.Range("D7").Value = 101
$wsAvance.Range("F7").Value = 8

$wsAvance.Range("C8").Formula = "=Daos!A105"
$wsAvance.Range("D8").Formula = "=Daos!C108"
$wsAvance.Range("F8").Value = 20

$wsAvance.Range("C9").Formula = "=Controllers!A94"
$wsAvance.Range("D9").Formula = "=Controllers!D97"
$wsAvance.Range("F9").Value = 40

$wsAvance.Range("C10").Value = 12
$wsAvance.Range("D10").Value = 12
$wsAvance.Range("F10").Value = 8

$wsAvance.Range("C11").Value = 83
$wsAvance.Range("D11").Value = 4
$wsAvance.Range("F11").Value = 120

# % de avance and Dias de avance columns, rows 4-11
for ($r = 4; $r -le 11; $r++) {
    $wsAvance.Range("E$r").Formula = "=D$r/C$r"
    $wsAvance.Range("G$r").Formula = "=F$r*E$r"
}

$detailRange = $wsAvance.Range("B4:G11")
$detailRange.Borders.LineStyle = 1
$wsAvance.Range("E4:E11").NumberFormat = "0%"
$wsAvance.Range("G4:G11").NumberFormat = "0"

# --- Totals row (12) -----------------------------------------------------
$wsAvance.Range("B12").Value = "Totales"
$wsAvance.Range("F12").Formula = "=SUM(F4:F11)"
$wsAvance.Range("G12").Formula = "=SUM(G4:G11)"
$wsAvance.Range("B12:G12").Borders.LineStyle = 1
$wsAvance.Range("G12").NumberFormat = "0"
$wsAvance.Range("B12").Font.Bold = $true
$wsAvance.Range("B12").Font.Size = 12
$wsAvance.Range("B12").HorizontalAlignment = -4152
$wsAvance.Range("B12").RowHeight = 15.75

# --- Avance Total (row 15) -------------------------------------------------
$wsAvance.Range("B15").Value = "Avance Total"
$wsAvance.Range("B15").Font.Size = 26
$wsAvance.Range("B15").RowHeight = 33.75

$wsAvance.Range("C15").Formula = "=G12/F12"
$wsAvance.Range("C15:D15").NumberFormat = "0.00%"
$wsAvance.Range("C15:D15").HorizontalAlignment = -4108
$wsAvance.Range("C15:D15").Merge()

$wsAvance.Range("D8").Select()

# ---------------------------------------------------------------------
# 3) Restore selections / active sheet state
# ---------------------------------------------------------------------
$wsDaos.Range("C54").Select()

$wsControllers.Activate()
$wsControllers.Range("D54").Select()
